$d = $word.ActiveDocument
$d.Content.Find.Execute("pratical use of A", $true, $false, $false, $false, $false,
                         $true, 1, $false, "practical use of A", 2)
